$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "Interest details" overflow rows (Applicant owns land,
# Permission obtained, Permission not obtained details). Everything below
# shifts up by three rows so that the previously row-68..87 content
# (Pre-application advice, Site details, Site Visit Details) becomes the
# new rows 65..84.
$ws.Range("A65:N67").EntireRow.Delete()

# Narrow columns D and E to their new widths (stored OOXML widths of 23
# and 31 respectively; Excel's ColumnWidth COM property reports ~0.83
# less than the stored width for this workbook's default font).
$ws.Columns.Item(4).ColumnWidth = 22.17
$ws.Columns.Item(5).ColumnWidth = 30.17
